$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2  (used throughout the Find.Execute calls below)

# --- 1) "-Oscar which is our Product owner (As in he manages the team)"
#        -> "-Oscar: Project manager"
$d.Content.Find.Execute(
    "-Oscar which is our Product owner (As in he manages the team)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-Oscar: Project manager", 2)

# --- 2) "-Camilo: " -> "-Camilo: Checker"
$d.Content.Find.Execute(
    "-Camilo: ", $true, $false, $false, $false, $false, $true, 1, $false,
    "-Camilo: Checker", 2)

# --- 3) "-Salva: " -> "-Salva: Writer"
$d.Content.Find.Execute(
    "-Salva: ", $true, $false, $false, $false, $false, $true, 1, $false,
    "-Salva: Writer", 2)

# --- 4) "current subject" -> "current subject."  (a period is appended)
$d.Content.Find.Execute(
    "current subject", $true, $false, $false, $false, $false, $true, 1, $false,
    "current subject.", 2)

# The original sentence had a trailing separate "." run right after the
# (now removed) _GoBack bookmark; collapse the resulting double period
# back down to a single one.
$d.Content.Find.Execute(
    "current subject..", $true, $false, $false, $false, $false, $true, 1, $false,
    "current subject.", 2)

# --- 5) Append two new paragraphs at the end of the document: an empty
#        paragraph followed by one describing the lends calendar feature.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$blankPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "With lends we added a calendar image for each month to make it easier to track the ongoing lends.X"

# --- 6) Move the _GoBack bookmark so it again sits at the very end of the
#        document content (right after the new last sentence). Placing a
#        bookmark collapsed exactly at the end-of-story position is
#        unreliable, so bookmark a throw-away trailing "X" character and
#        then delete it; the bookmark collapses cleanly to that position.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfText = $finalPara.Range.End - 2
$markRange = $d.Range($endOfText, $endOfText + 1)
$d.Bookmarks.Add("_GoBack", $markRange) | Out-Null
$markRange.Delete()
